$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 331.55554
$ws.Range("I9").Value = 330.57144
$ws.Range("J9").Value = 335
$ws.Range("K9").Value = 330.57144
$ws.Range("L9").Value = 335
$ws.Range("M9").Value = -161.57144
$ws.Range("H74").Value = 2640
$ws.Range("I74").Value = 2100
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2100
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -1164
$ws.Range("N74").Value = -4872
$ws.Range("H77").Value = 2640
$ws.Range("I77").Value = 2100
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 10500
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -5820
$ws.Range("N77").Value = -24360
$ws.Range("H137").Value = 1210.8
$ws.Range("I137").Value = 1327.8462
$ws.Range("J137").Value = 450
$ws.Range("K137").Value = 3983.5386
$ws.Range("L137").Value = 1350
$ws.Range("M137").Value = -1433.5386
$ws.Range("N137").Value = -6450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6748.7617
$ws.Range("I32").Value = 5334
$ws.Range("J32").Value = 11276
$ws.Range("K32").Value = 5334
$ws.Range("L32").Value = 11276
$ws.Range("M32").Value = -5047
$ws.Range("N32").Value = -11850

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7202.88
$ws.Range("I31").Value = 7152.3687
$ws.Range("J31").Value = 7362.8335
$ws.Range("K31").Value = 7152.3687
$ws.Range("L31").Value = 7362.8335
$ws.Range("M31").Value = -6857.3687
$ws.Range("N31").Value = -7952.8335
$ws.Range("H34").Value = 7202.88
$ws.Range("I34").Value = 7152.3687
$ws.Range("J34").Value = 7362.8335
$ws.Range("K34").Value = 7152.3687
$ws.Range("L34").Value = 7362.8335
$ws.Range("M34").Value = -6950.3687
$ws.Range("N34").Value = -7766.8335
$ws.Range("H86").Value = 4999.5
$ws.Range("I86").Value = 4999
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4999
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -3876
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4999.5
$ws.Range("I89").Value = 4999
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 24995
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -19379
$ws.Range("N89").Value = -36232
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H106").Value = 17671
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 17671
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 17671
$ws.Range("N106").Value = -20195
$ws.Range("H122").Value = 2166.3333
$ws.Range("I122").Value = 2166.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6498.999899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4048.999899999999
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 7634.25
$ws.Range("I132").Value = 7592.3335
$ws.Range("J132").Value = 7760
$ws.Range("K132").Value = 22777.0005
$ws.Range("L132").Value = 23280
$ws.Range("M132").Value = -20247.0005
$ws.Range("H134").Value = 1887.1666
$ws.Range("I134").Value = 1922.3636
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 5767.0908
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -3232.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 333474.66
$ws.Range("I4").Value = 212.5
$ws.Range("J4").Value = 999999
$ws.Range("K4").Value = 637.5
$ws.Range("L4").Value = 2999997
$ws.Range("M4").Value = -525.5
$ws.Range("N4").Value = -3000221
$ws.Range("H18").Value = 1807.8334
$ws.Range("I18").Value = 1807.8334
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 5423.5002
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -5254.5002
$ws.Range("H68").Value = 852
$ws.Range("I68").Value = 776.5
$ws.Range("J68").Value = 1003
$ws.Range("K68").Value = 2329.5
$ws.Range("L68").Value = 3009
$ws.Range("M68").Value = -1518.5
$ws.Range("H71").Value = 852
$ws.Range("I71").Value = 776.5
$ws.Range("J71").Value = 1003
$ws.Range("K71").Value = 6988.5
$ws.Range("L71").Value = 9027
$ws.Range("M71").Value = -2932.5
$ws.Range("H134").Value = 3739.8
$ws.Range("I134").Value = 3739.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11219.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6149.400000000001
$ws.Range("H136").Value = 3479.0833
$ws.Range("I136").Value = 3479.0833
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10437.2499
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5337.249899999999
$ws.Range("H139").Value = 1209.4
$ws.Range("I139").Value = 1209.4
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 3628.2
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 1511.8
$ws.Range("H140").Value = 808.2308
$ws.Range("I140").Value = 808.2308
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 2424.6924
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 2755.3076
$ws.Range("H141").Value = 3666.3333
$ws.Range("I141").Value = 1899.6
$ws.Range("J141").Value = 12500
$ws.Range("K141").Value = 5698.799999999999
$ws.Range("L141").Value = 37500
$ws.Range("M141").Value = -518.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4246.3335
$ws.Range("I102").Value = 4246.3335
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4246.3335
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2624.3335
$ws.Range("H122").Value = 25587.5
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 32950
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 98850
$ws.Range("M122").Value = -8050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 888.7
$ws.Range("I46").Value = 798.75
$ws.Range("J46").Value = 948.6667
$ws.Range("K46").Value = 798.75
$ws.Range("L46").Value = 948.6667
$ws.Range("M46").Value = -610.75
$ws.Range("N46").Value = -1324.6667
$ws.Range("H82").Value = 1431
$ws.Range("I82").Value = 1856.4286
$ws.Range("J82").Value = 686.5
$ws.Range("K82").Value = 1856.4286
$ws.Range("L82").Value = 686.5
$ws.Range("M82").Value = -1495.4286
$ws.Range("N82").Value = -1408.5
$ws.Range("H85").Value = 1431
$ws.Range("I85").Value = 1856.4286
$ws.Range("J85").Value = 686.5
$ws.Range("K85").Value = 1856.4286
$ws.Range("L85").Value = 686.5
$ws.Range("M85").Value = -608.4286
$ws.Range("N85").Value = -3182.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H132").Value = 4350.2856
$ws.Range("I132").Value = 4350.2856
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13050.8568
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -10520.8568
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
